$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to Text so numeric-looking strings
# (e.g. "224.71") are not auto-converted to numbers, matching the
# original inline-string cell contents.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.363.24'
$ws.Range('E2').Value = '  -1.21%  '
$ws.Range('D3').Value = '1.710.22'
$ws.Range('E3').Value = '  -1.11%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '224.71'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D6').Value = '0.5327'
$ws.Range('E6').Value = '  -1.96%  '
$ws.Range('E8').Value = '  -2.05%  '
$ws.Range('D9').Value = '0.06629'
$ws.Range('E9').Value = '  -0.56%  '
$ws.Range('D10').Value = '20.98'
$ws.Range('E10').Value = '  -3.73%  '
$ws.Range('D11').Value = '0.07630'
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = '4.556'
$ws.Range('E12').Value = '  -2.66%  '
$ws.Range('D13').Value = '1.729.58'
$ws.Range('E13').Value = '  +0.17%  '
$ws.Range('D14').Value = '1.946.02'
$ws.Range('E14').Value = '  -1.13%  '
$ws.Range('D15').Value = '0.5781'
$ws.Range('E15').Value = '  -2.76%  '
$ws.Range('D16').Value = '0.0₅8188'
$ws.Range('E16').Value = '  -2.31%  '
$ws.Range('D17').Value = '67.85'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = '27.381.78'
$ws.Range('E18').Value = '  -1.14%  '
$ws.Range('D19').Value = '218.57'
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').Value = '4.664'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('E22').Value = '  -3.52%  '
$ws.Range('D23').Value = '5.965'
$ws.Range('E23').Value = '  -3.70%  '
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '142.45'
$ws.Range('E25').Value = '  -3.25%  '
$ws.Range('D26').Value = '1.729'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = '0.1213'
$ws.Range('E27').Value = '  -2.69%  '
$ws.Range('E28').Value = '  -2.45%  '
$ws.Range('D29').Value = '16.25'
$ws.Range('E29').Value = '  -4.53%  '
$ws.Range('D30').Value = '0.05405'
$ws.Range('E30').Value = '  -4.48%  '
$ws.Range('D31').Value = '1.293'
$ws.Range('E31').Value = '  -1.32%  '
$ws.Range('D32').Value = '3.503'
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('D33').Value = '3.431'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('D34').Value = '1.649'
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('D35').Value = '2.877'
$ws.Range('E35').Value = '  +1.16%  '
$ws.Range('E36').Value = '  -2.39%  '
$ws.Range('E37').Value = '  -1.06%  '
$ws.Range('D38').Value = '0.5876'
$ws.Range('E38').Value = '  -1.57%  '
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = '5.860'
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('D41').Value = '1.048.18'
$ws.Range('E41').Value = '  +0.14%  '
$ws.Range('D42').Value = '1.004'
$ws.Range('E42').Value = '  +0.04%  '
$ws.Range('D43').Value = '0.8426'
$ws.Range('E43').Value = '  -1.96%  '
$ws.Range('D44').Value = '101.06'
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Value = '1.852.86'
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('E46').Value = '  +2.20%  '
$ws.Range('D47').Value = '58.06'
$ws.Range('D48').Value = '0.4517'
$ws.Range('E48').Value = '  +1.97%  '
$ws.Range('D49').Value = '1.003'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').Value = '8.110'
$ws.Range('E50').Value = '  -1.60%  '
$ws.Range('E51').Value = '  -1.79%  '

# Restore default (no explicit style) formatting so the cells match
# the original un-styled inline string cells.
$ws.Range("D2:E51").ClearFormats()
